$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 9 (diff hunk -1082,25 +1082,25)
$ws.Range("H9").Value = 7692409.5
$ws.Range("I9").Value = 12500062
$ws.Range("J9").Value = 166
$ws.Range("K9").Value = 12500062
$ws.Range("L9").Value = 166
$ws.Range("M9").Value = -12499893
$ws.Range("N9").Value = -504
# row 15 (diff hunk -1379,22 +1379,22)
$ws.Range("H15").Value = 1767.3944
$ws.Range("I15").Value = 1767.3944
$ws.Range("K15").Value = 5302.183199999999
$ws.Range("M15").Value = -5133.183199999999
# row 18 (diff hunk -1529,25 +1529,25)
$ws.Range("H18").Value = 3089595
$ws.Range("I18").Value = 6945889
$ws.Range("J18").Value = 4560
$ws.Range("K18").Value = 6945889
$ws.Range("L18").Value = 4560
$ws.Range("M18").Value = -6945605
$ws.Range("N18").Value = -5128
# row 40 (diff hunk -2631,20 +2631,23)
$ws.Range("H40").Value = 1625
$ws.Range("I40").Value = 1000
$ws.Range("K40").Value = 1000
$ws.Range("M40").Value = -825
# row 43 (diff hunk -2784,25 +2787,25)
$ws.Range("H43").Value = 1187
$ws.Range("J43").Value = 975
$ws.Range("L43").Value = 975
$ws.Range("N43").Value = -1113
# row 70 (diff hunk -4149,25 +4152,25)
$ws.Range("H70").Value = 1370.875
$ws.Range("I70").Value = 1302.6154
$ws.Range("J70").Value = 1666.6666
$ws.Range("K70").Value = 3907.8462
$ws.Range("L70").Value = 4999.9998
$ws.Range("M70").Value = -3637.8462
$ws.Range("N70").Value = -5539.9998
# row 73 (diff hunk -4302,25 +4305,25)
$ws.Range("H73").Value = 1370.875
$ws.Range("I73").Value = 1302.6154
$ws.Range("J73").Value = 1666.6666
$ws.Range("K73").Value = 3907.8462
$ws.Range("L73").Value = 4999.9998
$ws.Range("M73").Value = -2971.8462
$ws.Range("N73").Value = -6871.9998
# row 95 (diff hunk -5413,22 +5416,22)
$ws.Range("H95").Value = 27446.428
$ws.Range("J95").Value = 27446.428
$ws.Range("L95").Value = 27446.428
$ws.Range("N95").Value = -32938.428
# row 112 (diff hunk -6264,25 +6267,25)
$ws.Range("H112").Value = 1412.09
$ws.Range("J112").Value = 1458.2979
$ws.Range("L112").Value = 4374.893700000001
$ws.Range("N112").Value = -6590.893700000001
# row 127 (diff hunk -7020,25 +7023,25)
$ws.Range("H127").Value = 2018400.8
$ws.Range("I127").Value = 3154.4285
$ws.Range("J127").Value = 3585814.8
$ws.Range("K127").Value = 9463.2855
$ws.Range("L127").Value = 10757444.4
$ws.Range("M127").Value = -4503.2855
$ws.Range("N127").Value = -10767364.4
# row 129 (diff hunk -7121,25 +7124,25)
$ws.Range("H129").Value = 4546460
$ws.Range("I129").Value = 62501670
$ws.Range("J129").Value = 953.549
$ws.Range("K129").Value = 187505010
$ws.Range("L129").Value = 2860.647
$ws.Range("M129").Value = -187500010
$ws.Range("N129").Value = -12860.647
# row 138 (diff hunk -7574,25 +7577,25)
$ws.Range("H138").Value = 7056.7446
$ws.Range("I138").Value = 3234.7727
$ws.Range("J138").Value = 8293.264999999999
$ws.Range("K138").Value = 9704.3181
$ws.Range("L138").Value = 24879.795
$ws.Range("M138").Value = -4564.3181
$ws.Range("N138").Value = -35159.795
# row 140 (diff hunk -7675,22 +7678,22)
$ws.Range("H140").Value = 35000
$ws.Range("J140").Value = 35000
$ws.Range("L140").Value = 35000
$ws.Range("N140").Value = -45360

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 3 (diff hunk -7922,25 +7925,25)
$ws.Range("H3").Value = 6094.2856
$ws.Range("I3").Value = 2732
$ws.Range("J3").Value = 14500
$ws.Range("K3").Value = 2732
$ws.Range("L3").Value = 14500
$ws.Range("M3").Value = -2617
$ws.Range("N3").Value = -14730
# row 32 (diff hunk -9355,22 +9358,22)
$ws.Range("H32").Value = 21153.74
$ws.Range("I32").Value = 17557.6
$ws.Range("K32").Value = 17557.6
$ws.Range("M32").Value = -17270.6
# row 74 (diff hunk -11407,25 +11410,25)
$ws.Range("H74").Value = 2387.8948
$ws.Range("I74").Value = 1740.7858
$ws.Range("J74").Value = 4199.8
$ws.Range("K74").Value = 1740.7858
$ws.Range("L74").Value = 4199.8
$ws.Range("M74").Value = -866.7858000000001
$ws.Range("N74").Value = -5947.8
# row 77 (diff hunk -11554,25 +11557,25)
$ws.Range("H77").Value = 2387.8948
$ws.Range("I77").Value = 1740.7858
$ws.Range("J77").Value = 4199.8
$ws.Range("K77").Value = 8703.929
$ws.Range("L77").Value = 20999
$ws.Range("M77").Value = -4335.929
$ws.Range("N77").Value = -29735

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 35 (diff hunk -16453,22 +16456,22)
$ws.Range("H35").Value = 22916
$ws.Range("J35").Value = 22916
$ws.Range("L35").Value = 22916
$ws.Range("N35").Value = -23536

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31 (diff hunk -23214,25 +23217,25)
$ws.Range("H31").Value = 3337376.5
$ws.Range("I31").Value = 4350304.5
$ws.Range("J31").Value = 9184.857
$ws.Range("K31").Value = 4350304.5
$ws.Range("L31").Value = 9184.857
$ws.Range("M31").Value = -4350009.5
$ws.Range("N31").Value = -9774.857
# row 34 (diff hunk -23370,25 +23373,25)
$ws.Range("H34").Value = 3337376.5
$ws.Range("I34").Value = 4350304.5
$ws.Range("J34").Value = 9184.857
$ws.Range("K34").Value = 4350304.5
$ws.Range("L34").Value = 9184.857
$ws.Range("M34").Value = -4350102.5
$ws.Range("N34").Value = -9588.857
# row 122 (diff hunk -27670,25 +27673,25)
$ws.Range("H122").Value = 2843.4138
$ws.Range("I122").Value = 2531.3333
$ws.Range("J122").Value = 3662.625
$ws.Range("K122").Value = 7593.999899999999
$ws.Range("L122").Value = 10987.875
$ws.Range("M122").Value = -5143.999899999999
$ws.Range("N122").Value = -15887.875
# row 123 (diff hunk -27722,19 +27725,22)
$ws.Range("H123").Value = 30993.334
$ws.Range("J123").Value = 30993.334
$ws.Range("L123").Value = 30993.334
$ws.Range("N123").Value = -40793.334

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 5 (diff hunk -28900,25 +28906,25)
$ws.Range("H5").Value = 926.7143
$ws.Range("I5").Value = 583.4286
$ws.Range("J5").Value = 2986.4285
$ws.Range("K5").Value = 1750.2858
$ws.Range("L5").Value = 8959.2855
$ws.Range("M5").Value = -1638.2858
$ws.Range("N5").Value = -9183.2855
# row 87 (diff hunk -33083,22 +33089,22)
$ws.Range("H87").Value = 12971.571
$ws.Range("I87").Value = 10700.25
$ws.Range("K87").Value = 32100.75
$ws.Range("M87").Value = -30852.75
# row 90 (diff hunk -33239,22 +33245,22)
$ws.Range("H90").Value = 12971.571
$ws.Range("I90").Value = 10700.25
$ws.Range("K90").Value = 96302.25
$ws.Range("M90").Value = -90062.25
# row 107 (diff hunk -34108,25 +34114,25)
$ws.Range("H107").Value = 1136.1428
$ws.Range("I107").Value = 361.875
$ws.Range("J107").Value = 1612.6154
$ws.Range("K107").Value = 1085.625
$ws.Range("L107").Value = 4837.8462
$ws.Range("M107").Value = 834.375
$ws.Range("N107").Value = -8677.8462
# row 131 (diff hunk -35347,25 +35353,25)
$ws.Range("H131").Value = 1710.1136
$ws.Range("J131").Value = 1279.1464
$ws.Range("L131").Value = 3837.4392
$ws.Range("N131").Value = -13917.4392
# row 133 (diff hunk -35451,25 +35457,25)
$ws.Range("H133").Value = 5835
$ws.Range("J133").Value = 5000
$ws.Range("L133").Value = 15000
$ws.Range("N133").Value = -25120
# row 134 (diff hunk -35503,25 +35509,22)
$ws.Range("H134").Value = 1033
$ws.Range("I134").Value = 1033
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3099
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 1971
$ws.Range("N134").ClearContents()
# row 135 (diff hunk -35555,25 +35558,25)
$ws.Range("H135").Value = 926.7143
$ws.Range("I135").Value = 583.4286
$ws.Range("J135").Value = 2986.4285
$ws.Range("K135").Value = 5250.8574
$ws.Range("L135").Value = 26877.8565
$ws.Range("M135").Value = -2715.8574
$ws.Range("N135").Value = -31947.8565
# row 137 (diff hunk -35659,25 +35662,22)
$ws.Range("H137").Value = 3025
$ws.Range("I137").Value = 3025
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 9075
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -3975
$ws.Range("N137").ClearContents()
# row 138 (diff hunk -35711,25 +35711,22)
$ws.Range("H138").Value = 1239.8889
$ws.Range("I138").Value = 1239.8889
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 3719.6667
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 1420.3333
$ws.Range("N138").ClearContents()
# row 139 (diff hunk -35763,25 +35760,25)
$ws.Range("H139").Value = 10872871
$ws.Range("I139").Value = 11366320
$ws.Range("J139").Value = 17000
$ws.Range("K139").Value = 34098960
$ws.Range("L139").Value = 51000
$ws.Range("M139").Value = -34093820
$ws.Range("N139").Value = -61280

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 80 (diff hunk -39841,25 +39838,25)
$ws.Range("H80").Value = 3238.4092
$ws.Range("I80").Value = 3262.5
$ws.Range("J80").Value = 3174.1667
$ws.Range("K80").Value = 3262.5
$ws.Range("L80").Value = 3174.1667
$ws.Range("M80").Value = -2264.5
$ws.Range("N80").Value = -5170.1667
# row 83 (diff hunk -39991,25 +39988,25)
$ws.Range("H83").Value = 3238.4092
$ws.Range("I83").Value = 3262.5
$ws.Range("J83").Value = 3174.1667
$ws.Range("K83").Value = 16312.5
$ws.Range("L83").Value = 15870.8335
$ws.Range("M83").Value = -11320.5
$ws.Range("N83").Value = -25854.8335
# row 104 (diff hunk -41029,22 +41026,22)
$ws.Range("H104").Value = 29800
$ws.Range("J104").Value = 29800
$ws.Range("L104").Value = 29800
$ws.Range("N104").Value = -36788

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 82 (diff hunk -46860,25 +46857,25)
$ws.Range("H82").Value = 2269.8125
$ws.Range("I82").Value = 1664.875
$ws.Range("J82").Value = 2874.75
$ws.Range("K82").Value = 1664.875
$ws.Range("L82").Value = 2874.75
$ws.Range("M82").Value = -1303.875
$ws.Range("N82").Value = -3596.75
# row 85 (diff hunk -47010,25 +47007,25)
$ws.Range("H85").Value = 2269.8125
$ws.Range("I85").Value = 1664.875
$ws.Range("J85").Value = 2874.75
$ws.Range("K85").Value = 1664.875
$ws.Range("L85").Value = 2874.75
$ws.Range("M85").Value = -416.875
$ws.Range("N85").Value = -5370.75
# row 122 (diff hunk -48793,25 +48790,25)
$ws.Range("H122").Value = 3929.2
$ws.Range("I122").Value = 2874.5
$ws.Range("J122").Value = 4632.3335
$ws.Range("K122").Value = 8623.5
$ws.Range("L122").Value = 13897.0005
$ws.Range("M122").Value = -6173.5
$ws.Range("N122").Value = -18797.0005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 58 (diff hunk -52617,22 +52614,22)
$ws.Range("H58").Value = 13769.23
$ws.Range("I58").Value = 11000
$ws.Range("K58").Value = 11000
$ws.Range("M58").Value = -10692
# row 75 (diff hunk -53459,22 +53456,22)
$ws.Range("H75").Value = 35633.332
$ws.Range("J75").Value = 35633.332
$ws.Range("L75").Value = 35633.332
$ws.Range("N75").Value = -37505.332
# row 78 (diff hunk -53606,22 +53603,22)
$ws.Range("H78").Value = 35633.332
$ws.Range("J78").Value = 35633.332
$ws.Range("L78").Value = 106899.996
$ws.Range("N78").Value = -116259.996
# row 138 (diff hunk -56534,22 +56531,19)
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
